# Tratamento inicial de dados
$wb = $excel.ActiveWorkbook

# --- Rename existing sheets / reorganize tabs ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Colunas e Relação com Vazio"

# Keep a handle to the original "Plan2" sheet (it will shift to the 3rd
# position once the new sheet is inserted after $ws1).
$wsPlan2 = $wb.Worksheets.Item(2)

# --- Insert the new sheet "Linhas com Satisfação 99" between the two ---
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Linhas com Satisfação 99"

# --- Populate the new sheet with its header + 82 data rows ---
$newSheet.Range("A1").Value = "Total 82"
$newSheet.Cells.Item(2, 1).Value = 460371
$newSheet.Cells.Item(3, 1).Value = 393512
$newSheet.Cells.Item(4, 1).Value = 617112
$newSheet.Cells.Item(5, 1).Value = 557223
$newSheet.Cells.Item(6, 1).Value = 486657
$newSheet.Cells.Item(7, 1).Value = 664283
$newSheet.Cells.Item(8, 1).Value = 504864
$newSheet.Cells.Item(9, 1).Value = 373736
$newSheet.Cells.Item(10, 1).Value = 655270
$newSheet.Cells.Item(11, 1).Value = 214162
$newSheet.Cells.Item(12, 1).Value = 491249
$newSheet.Cells.Item(13, 1).Value = 524992
$newSheet.Cells.Item(14, 1).Value = 928539
$newSheet.Cells.Item(15, 1).Value = 930916
$newSheet.Cells.Item(16, 1).Value = 72862
$newSheet.Cells.Item(17, 1).Value = 91089
$newSheet.Cells.Item(18, 1).Value = 113944
$newSheet.Cells.Item(19, 1).Value = 195564
$newSheet.Cells.Item(20, 1).Value = 573385
$newSheet.Cells.Item(21, 1).Value = 624534
$newSheet.Cells.Item(22, 1).Value = 867398
$newSheet.Cells.Item(23, 1).Value = 30794
$newSheet.Cells.Item(24, 1).Value = 625115
$newSheet.Cells.Item(25, 1).Value = 630294
$newSheet.Cells.Item(26, 1).Value = 630440
$newSheet.Cells.Item(27, 1).Value = 1012388
$newSheet.Cells.Item(28, 1).Value = 101929
$newSheet.Cells.Item(29, 1).Value = 669422
$newSheet.Cells.Item(30, 1).Value = 327571
$newSheet.Cells.Item(31, 1).Value = 337927
$newSheet.Cells.Item(32, 1).Value = 698290
$newSheet.Cells.Item(33, 1).Value = 426532
$newSheet.Cells.Item(34, 1).Value = 442103
$newSheet.Cells.Item(35, 1).Value = 518204
$newSheet.Cells.Item(36, 1).Value = 557243
$newSheet.Cells.Item(37, 1).Value = 726529
$newSheet.Cells.Item(38, 1).Value = 828295
$newSheet.Cells.Item(39, 1).Value = 926152
$newSheet.Cells.Item(40, 1).Value = 931518
$newSheet.Cells.Item(41, 1).Value = 740091
$newSheet.Cells.Item(42, 1).Value = 771606
$newSheet.Cells.Item(43, 1).Value = 406214
$newSheet.Cells.Item(44, 1).Value = 740978
$newSheet.Cells.Item(45, 1).Value = 464844
$newSheet.Cells.Item(46, 1).Value = 541195
$newSheet.Cells.Item(47, 1).Value = 288917
$newSheet.Cells.Item(48, 1).Value = 358796
$newSheet.Cells.Item(49, 1).Value = 780204
$newSheet.Cells.Item(50, 1).Value = 756742
$newSheet.Cells.Item(51, 1).Value = 584161
$newSheet.Cells.Item(52, 1).Value = 836769
$newSheet.Cells.Item(53, 1).Value = 812304
$newSheet.Cells.Item(54, 1).Value = 859676
$newSheet.Cells.Item(55, 1).Value = 836758
$newSheet.Cells.Item(56, 1).Value = 625736
$newSheet.Cells.Item(57, 1).Value = 585286
$newSheet.Cells.Item(58, 1).Value = 862294
$newSheet.Cells.Item(59, 1).Value = 407542
$newSheet.Cells.Item(60, 1).Value = 13886
$newSheet.Cells.Item(61, 1).Value = 177178
$newSheet.Cells.Item(62, 1).Value = 910504
$newSheet.Cells.Item(63, 1).Value = 467919
$newSheet.Cells.Item(64, 1).Value = 186744
$newSheet.Cells.Item(65, 1).Value = 248022
$newSheet.Cells.Item(66, 1).Value = 310568
$newSheet.Cells.Item(67, 1).Value = 419862
$newSheet.Cells.Item(68, 1).Value = 578470
$newSheet.Cells.Item(69, 1).Value = 21239312
$newSheet.Cells.Item(70, 1).Value = 20222726
$newSheet.Cells.Item(71, 1).Value = 10876194
$newSheet.Cells.Item(72, 1).Value = 10532189
$newSheet.Cells.Item(73, 1).Value = 20245604
$newSheet.Cells.Item(74, 1).Value = 11149425
$newSheet.Cells.Item(75, 1).Value = 10419702
$newSheet.Cells.Item(76, 1).Value = 10161268
$newSheet.Cells.Item(77, 1).Value = 10226764
$newSheet.Cells.Item(78, 1).Value = 21140976
$newSheet.Cells.Item(79, 1).Value = 10747120
$newSheet.Cells.Item(80, 1).Value = 10821051
$newSheet.Cells.Item(81, 1).Value = 10918535
$newSheet.Cells.Item(82, 1).Value = 20544139
$newSheet.Cells.Item(83, 1).Value = 20001614

# Leave the new sheet's selection on A2, matching the target view state.
$null = $newSheet.Range("A2").Select()

# --- Adjust the view state of the first sheet (no longer the active tab) ---
$ws1.Activate()
$null = $ws1.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

# --- Re-activate the new sheet so it becomes the workbook's active tab ---
$newSheet.Activate()
$null = $newSheet.Range("A2").Select()
